# Resort the "Periodo Mora" / "Valor Mora" table (rows 16-107 on Hoja1)
# from descending (newest first) to ascending (oldest first) chronological
# order, carrying each row's "Valor Mora" along with its period so the
# figures stay attached to the correct period after the resort.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodosYValores = @(
    @("1607", 61840),
    @("1608", 61840),
    @("1609", 61840),
    @("1610", 61840),
    @("1611", 61840),
    @("1612", 61840),
    @("1701", 61840),
    @("1702", 61840),
    @("1703", 61840),
    @("1704", 61840),
    @("1705", 61840),
    @("1706", 61840),
    @("1707", 61840),
    @("1708", 61840),
    @("1709", 61840),
    @("1710", 61840),
    @("1711", 61840),
    @("1712", 61840),
    @("1801", 61840),
    @("1802", 61840),
    @("1803", 61840),
    @("1804", 61840),
    @("1805", 61840),
    @("1806", 61840),
    @("1807", 61840),
    @("1808", 61840),
    @("1809", 61840),
    @("1810", 61840),
    @("1811", 61840),
    @("1812", 61840),
    @("1901", 61840),
    @("1902", 61840),
    @("1903", 61840),
    @("1904", 61840),
    @("1905", 61840),
    @("1906", 61840),
    @("1907", 61840),
    @("1908", 61840),
    @("1909", 61840),
    @("1910", 61840),
    @("1911", 61840),
    @("1912", 61840),
    @("2001", 61840),
    @("2002", 61840),
    @("2003", 61840),
    @("2004", 61840),
    @("2005", 61840),
    @("2006", 61840),
    @("2007", 61840),
    @("2008", 61840),
    @("2009", 61840),
    @("2010", 61840),
    @("2011", 61840),
    @("2012", 61840),
    @("2101", 61840),
    @("2102", 61840),
    @("2103", 61840),
    @("2104", 61840),
    @("2105", 61840),
    @("2106", 61840),
    @("2107", 61840),
    @("2108", 61840),
    @("2109", 61840),
    @("2110", 61840),
    @("2111", 61840),
    @("2112", 61840),
    @("2201", 61840),
    @("2202", 61840),
    @("2203", 61840),
    @("2204", 61840),
    @("2205", 61840),
    @("2206", 61840),
    @("2207", 61840),
    @("2208", 61840),
    @("2209", 61840),
    @("2210", 61840),
    @("2211", 61840),
    @("2212", 61840),
    @("2301", 61840),
    @("2302", 61840),
    @("2303", 61840),
    @("2304", 61840),
    @("2305", 61840),
    @("2306", 61840),
    @("2307", 61840),
    @("2308", 61840),
    @("2309", 61840),
    @("2310", 61840),
    @("2311", 61840),
    @("2312", 61840),
    @("2401", 61840),
    @("2402", 37104),
)

$firstRow = 16
$row = $firstRow
foreach ($item in $periodosYValores) {
    $ws.Cells.Item($row, 5).Value = $item[0]
    $ws.Cells.Item($row, 6).Value = $item[1]
    $row++
}

